# Fix [s] vs [q] in Queue implementation
#
# Slide 16 ("Queue implementation, as list") has an OCaml code listing in the
# "Content Placeholder 2" shape. The module uses `q` as the parameter name for
# the queue everywhere except four spots that were mistakenly left as `s`:
#
#   let is_empty s = (s = [])        ->   let is_empty q = (q = [])
#   let enq x s = s @ [x]             ->   let enq x q = q @ [x]
#
# Each occurrence is replaced in place (same length, "s" -> "q") so the
# existing run/formatting (syntax-highlight colors, Courier font, etc.) is
# left untouched.

$p   = $ppt.ActivePresentation
$s16 = $p.Slides.Item(16)
$shp = $s16.Shapes.Item(2)
$tr  = $shp.TextFrame.TextRange

function Replace-NthOccurrence {
    param(
        $TextRange,
        [string]$Search,
        [string]$Replace,
        [int]$Occurrence = 1
    )

    $full = $TextRange.Text
    $searchStart = 0
    $found = -1

    for ($k = 0; $k -lt $Occurrence; $k++) {
        $found = $full.IndexOf($Search, $searchStart)
        if ($found -lt 0) { break }
        $searchStart = $found + 1
    }

    if ($found -ge 0) {
        # TextRange.Characters is 1-based.
        $sub = $TextRange.Characters($found + 1, $Search.Length)
        $sub.Text = $Replace
    }
}

# "let is_empty s = (s = [])"
Replace-NthOccurrence $tr " s " " q " 1
Replace-NthOccurrence $tr "s " "q " 1

# "let enq x s = s @ [x]"
Replace-NthOccurrence $tr " x s " " x q " 1
Replace-NthOccurrence $tr " s " " q " 1

Write-Output $tr.Text
